$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The log table (Date / Time start / Time end / Total time) gains two more
# shifts: 2024-11-10 (rows 7:8) and 2024-11-17 (rows 9:10), each following the
# existing pattern of a merged Date cell spanning two rows plus a start/end
# time pair per row. Merge the new Date cells first (like the earlier pairs)
# before touching their formatting/values, then copy the number formats down
# from the most recent existing pair so the new cells pick up the same
# (already-defined) date / time styles instead of manufacturing new ones.
# ---------------------------------------------------------------------------

$ws.Range("F7:F8").MergeCells = $true
$ws.Range("F9:F10").MergeCells = $true

$ws.Range("F5:F6").Copy() | Out-Null
$ws.Range("F7:F8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F9:F10").PasteSpecial(-4122) | Out-Null

$ws.Range("G4:I4").Copy() | Out-Null
$ws.Range("G7:I10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Shift on 2024-11-10 (45606): 18:30-20:00 and 21:00-22:00 -------------
$ws.Range("F7").Value = 45606
$ws.Range("G7").Value = 0.77083333333333337
$ws.Range("H7").Value = 0.83333333333333337
$ws.Range("I7").Formula = "=H7-G7"

$ws.Range("G8").Value = 0.875
$ws.Range("H8").Value = 0.91666666666666663
$ws.Range("I8").Formula = "=H8-G8"

# --- Shift on 2024-11-17 (45613): 17:30-19:00 and 20:00-22:30 -------------
$ws.Range("F9").Value = 45613
$ws.Range("G9").Value = 0.72916666666666663
$ws.Range("H9").Value = 0.79166666666666663
$ws.Range("I9").Formula = "=H9-G9"

$ws.Range("G10").Value = 0.83333333333333337
$ws.Range("H10").Value = 0.9375
$ws.Range("I10").Formula = "=H10-G10"

# C3/C4/C5 are formulas (SUM of total time, hourly payout, approx payout) and
# recalculate automatically from the new rows above.

# --- Selection cursor moved to P11 -----------------------------------------
$ws.Range("P11").Select() | Out-Null
